# Update the number of pieces of Cost4 (and Cost5) on the "Max N" sheet,
# and make "Max N" the active/selected tab with I2 selected.

$wb = $excel.ActiveWorkbook

$wsMaxN = $wb.Worksheets.Item("Max N")

# Update the piece counts: Cost4 (column E) 12 -> 10, Cost5 (column F) 10 -> 9
$wsMaxN.Range("E2").Value = 10
$wsMaxN.Range("F2").Value = 9

# Make "Max N" the active sheet/tab, with I2 as the selected cell.
$wsMaxN.Activate()
$wsMaxN.Range("I2").Select()
